# Apply the "additional scraping" update:
#  1. Insert a new "Player Info" sheet as the first sheet in the workbook,
#     containing ID / NAME / BATTING_HAND / BOWL_STYLE for player 4705.
#  2. On the existing "ODI Batting" and "ODI Bowling" sheets, rename the
#     MATCH_CARD_LINK column to MATCH_CODE and replace the full scorecard
#     URL values with just the numeric match code extracted from the URL.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Player Info" worksheet before the current first sheet
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Capture the existing bold/bordered header style BEFORE inserting the new
# sheet (reading range properties after Worksheets.Add can return stale data
# in this runtime, so grab everything we need up front).
$headerStyle = $battingSheet.Range("A1")
$hdrBold = $headerStyle.Font.Bold()
$hdrHAlign = $headerStyle.HorizontalAlignment()
$hdrVAlign = $headerStyle.VerticalAlignment()
$hdrBorder = $headerStyle.Borders.LineStyle()

$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold/bordered header style used on the other sheets
$hdrRange = $playerInfo.Range("A1:D1")
$hdrRange.Font.Bold = $hdrBold
$hdrRange.HorizontalAlignment = $hdrHAlign
$hdrRange.VerticalAlignment = $hdrVAlign
$hdrRange.Borders.LineStyle = $hdrBorder

# Data row - keep ID as text (matches the source inlineStr representation)
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4705"
$playerInfo.Range("B2").Value = "Jason Paul Behrendorff"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Fast Medium"

# ---------------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and keep only the match code
# ---------------------------------------------------------------------
function Update-MatchCodeColumn($ws, [int]$col, [int]$lastRow) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value = "MATCH_CODE"

    $dataRange = $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item($lastRow, $col))
    $dataRange.NumberFormat = "@"

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $link = $cell.Value()
        $code = ($link -split "MatchCode=")[1]
        $cell.Value = $code
    }
}

$odiBatting = $wb.Worksheets.Item("ODI Batting")
Update-MatchCodeColumn $odiBatting 4 13

$odiBowling = $wb.Worksheets.Item("ODI Bowling")
Update-MatchCodeColumn $odiBowling 2 13
